$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update mapped column names (B) that changed to the new DB naming convention
$ws.Range("B4").Value = "occupant_user_id"
$ws.Range("B5").Value = "substitute_user_id"
$ws.Range("B12").Value = "referente ao id func"
$ws.Range("B15").Value = "referente ao id func"

# Adjust column widths to better fit the updated (longer) content in column B
$ws.Columns.Item(1).ColumnWidth = 21.666666666666668
$ws.Columns.Item(2).ColumnWidth = 28.833333333333332

# Update the active selection as left by the author after editing
$ws.Range("B10").Select()
